$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.138.16"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +5.10%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.612.29"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +5.26%  "

# Row 4
$ws.Range("E4").Value = "  -0.24%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.23%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "190.24"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.30%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.643"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.11%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.598.97"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.08%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.11%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.177"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.86%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.660"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.89%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "58.36"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.77%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000289"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.19%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.85"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.85%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.180.46"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.72%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.55"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.04%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.603.10"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.77%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.069.63"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.11%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.56"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.08%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.121"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.48%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.05"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.78%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "491.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.61%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.32"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +12.31%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.40"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.34%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.46"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.22%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "90.70"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.83%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.11"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.52%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.08"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.10%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.38"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.97%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.81"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.70%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.58"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.10%  "

# Row 32
$ws.Range("B32").Value = "Cosmos"
$ws.Range("C32").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.32"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.54%  "

# Row 33
$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "623.71"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.00%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.118"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.10%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "65.29"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.21%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0₃0823"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.26%  "

# Row 37
$ws.Range("B37").Value = "InjectiveProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "38.27"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.52%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.403"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.80%  "

# Row 39
$ws.Range("B39").Value = "Dai"
$ws.Range("C39").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.13%  "

# Row 40
$ws.Range("E40").Value = "  -1.70%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.60"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.04%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.333.46"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.67%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.18"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +10.35%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0450"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.07%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.71"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.00%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.29"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.12%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.138"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.86%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.09"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.61%  "

# Row 49
$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.72"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.65%  "

# Row 50
$ws.Range("B50").Value = "LidoDAOToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.30"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.04%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.998"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.31%  "
